$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.045.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.513.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("E9").Value = "  +6.96%  "
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.121.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.135"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.033.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.521.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.540"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.894"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.797.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "340.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.844"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
